# Replace the "StatQuery" cell (C2, ParticipantsTab row) with the updated
# stats query: restructured joins (df_study/df_participant/df_sample/
# df_sequencing_file/df_pathology_file aliases), Files now counts both
# sequencing and pathology files, and an added sex_at_birth = 'Female' filter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
SELECT
    COUNT(DISTINCT std.study_ID) AS "Studies",
    COUNT(DISTINCT prt.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    (COUNT(DISTINCT seq.id) + COUNT(DISTINCT paf.id)) AS "Files"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN 
    df_pathology_file paf ON smp.id = paf."sample.id"
WHERE 
    std.study_ID = 'phs002430' 
    AND prt.race = 'Asian' 
    AND prt.sex_at_birth = 'Female';
'@

$ws.Range("C2").Value = $newQuery

# Mirror the author's scroll/selection change: the sheet was left scrolled
# to the top with C2 (the cell just edited) as the active cell, rather than
# the previous C6 selection.
[void]$ws.Range("C2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 3
